$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'26.513.97"
$ws.Cells.Item(2, 5).Value = "  +0.90%  "

$ws.Cells.Item(3, 4).Value = "'1.728.30"

$ws.Cells.Item(4, 5).Value = "  +0.04%  "

$ws.Cells.Item(5, 4).Value = "'245.25"
$ws.Cells.Item(5, 5).Value = "  +2.14%  "

$ws.Cells.Item(6, 5).Value = "  +0.01%  "

$ws.Cells.Item(7, 4).Value = "'0.4811"
$ws.Cells.Item(7, 5).Value = "  +1.72%  "

$ws.Cells.Item(8, 4).Value = "'0.2669"
$ws.Cells.Item(8, 5).Value = "  +1.66%  "

$ws.Cells.Item(9, 4).Value = "'0.06194"

$ws.Cells.Item(10, 4).Value = "'1.731.40"

$ws.Cells.Item(11, 4).Value = "'0.07183"
$ws.Cells.Item(11, 5).Value = "  +1.81%  "

$ws.Cells.Item(12, 5).Value = "  +0.69%  "

$ws.Cells.Item(13, 4).Value = "'0.6098"
$ws.Cells.Item(13, 5).Value = "  +1.99%  "

$ws.Cells.Item(14, 4).Value = "'4.528"
$ws.Cells.Item(14, 5).Value = "  +2.12%  "

$ws.Cells.Item(15, 4).Value = "'77.16"
$ws.Cells.Item(15, 5).Value = "  +1.29%  "

$ws.Cells.Item(16, 4).Value = "'0.9998"
$ws.Cells.Item(16, 5).Value = "  +0.00%  "

$ws.Cells.Item(17, 4).Value = "'26.512.58"
$ws.Cells.Item(17, 5).Value = "  +0.87%  "

$ws.Cells.Item(18, 4).Value = "'0.9998"
$ws.Cells.Item(18, 5).Value = "  +0.01%  "

$ws.Cells.Item(19, 5).Value = "  +1.85%  "

$ws.Cells.Item(20, 4).Value = "'11.53"
$ws.Cells.Item(20, 5).Value = "  +0.00%  "

$ws.Cells.Item(21, 4).Value = "'1.953.88"
$ws.Cells.Item(21, 5).Value = "  +0.93%  "

$ws.Cells.Item(22, 4).Value = "'4.524"
$ws.Cells.Item(22, 5).Value = "  -0.22%  "

$ws.Cells.Item(23, 4).Value = "'8.813"
$ws.Cells.Item(23, 5).Value = "  +1.19%  "

$ws.Cells.Item(24, 4).Value = "'5.253"
$ws.Cells.Item(24, 5).Value = "  +0.05%  "

$ws.Cells.Item(25, 4).Value = "'137.05"
$ws.Cells.Item(25, 5).Value = "  +1.42%  "

$ws.Cells.Item(26, 5).Value = "  +0.89%  "

$ws.Cells.Item(27, 4).Value = "'1.773"
$ws.Cells.Item(27, 5).Value = "  +0.17%  "

$ws.Cells.Item(28, 4).Value = "'1.406"
$ws.Cells.Item(28, 5).Value = "  +0.58%  "

$ws.Cells.Item(29, 4).Value = "'107.42"
$ws.Cells.Item(29, 5).Value = "  +0.59%  "

$ws.Cells.Item(30, 4).Value = "'3.974"
$ws.Cells.Item(30, 5).Value = "  +0.72%  "

$ws.Cells.Item(31, 4).Value = "'0.08031"
$ws.Cells.Item(31, 5).Value = "  +3.07%  "

$ws.Cells.Item(32, 4).Value = "'3.691"
$ws.Cells.Item(32, 5).Value = "  -0.01%  "

$ws.Cells.Item(33, 4).Value = "'0.04515"
$ws.Cells.Item(33, 5).Value = "  +0.23%  "

$ws.Cells.Item(34, 5).Value = "  +0.08%  "

$ws.Cells.Item(35, 4).Value = "'0.9996"
$ws.Cells.Item(35, 5).Value = "  +2.18%  "

$ws.Cells.Item(36, 4).Value = "'0.6250"
$ws.Cells.Item(36, 5).Value = "  +0.74%  "

$ws.Cells.Item(37, 4).Value = "'2.079"
$ws.Cells.Item(37, 5).Value = "  +7.78%  "

$ws.Cells.Item(38, 4).Value = "'0.9080"
$ws.Cells.Item(38, 5).Value = "  -2.57%  "

$ws.Cells.Item(39, 5).Value = "  -2.67%  "

$ws.Cells.Item(40, 4).Value = "'1.001"
$ws.Cells.Item(40, 5).Value = "  +0.10%  "

$ws.Cells.Item(41, 4).Value = "'0.01505"
$ws.Cells.Item(41, 5).Value = "  +1.63%  "

$ws.Cells.Item(42, 4).Value = "'102.35"
$ws.Cells.Item(42, 5).Value = "  -10.96%  "

$ws.Cells.Item(43, 4).Value = "'5.569"
$ws.Cells.Item(43, 5).Value = "  -0.87%  "

$ws.Cells.Item(44, 4).Value = "'0.3876"
$ws.Cells.Item(44, 5).Value = "  +1.31%  "

$ws.Cells.Item(45, 4).Value = "'6.960"
$ws.Cells.Item(45, 5).Value = "  +9.85%  "

$ws.Cells.Item(46, 5).Value = "  +0.26%  "

$ws.Cells.Item(47, 4).Value = "'0.05366"
$ws.Cells.Item(47, 5).Value = "  +2.00%  "

$ws.Cells.Item(48, 5).Value = "  +0.38%  "

$ws.Cells.Item(49, 4).Value = "'7.797"
$ws.Cells.Item(49, 5).Value = "  +0.31%  "

$ws.Cells.Item(50, 5).Value = "  +2.86%  "

$ws.Cells.Item(51, 2).Value = "Decentraland"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(51, 4).Value = "'0.3396"
$ws.Cells.Item(51, 5).Value = "  +0.36%  "
